$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the hidden "_GoBack" bookmark from the very first paragraph
#    of the document down to the (empty) ListParagraph that follows
#    the "Corbin: ..." bullet item.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r" -and $p.Style.NameLocal -eq "List Paragraph") {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $d.Bookmarks.Add("_GoBack", $targetPara.Range)
}

# ------------------------------------------------------------------
# 2) Append the new "CORBIN UPDATE" section after the existing
#    "Visual Studio Code" bullet at the end of the document.
# ------------------------------------------------------------------
$endPos = $d.Content.End - 1
$insertionRange = $d.Range($endPos, $endPos)

$newContent = '<w:p><w:pPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:highlight w:val="magenta"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="magenta"/></w:rPr><w:t xml:space="preserve">CORBIN UPDATE </w:t></w:r></w:p><w:p><w:r><w:t>Tools for application development:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Eclipse IDE</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Java SE 8</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>JavaFX</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Java </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>JDK</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Gluon Scene Viewer</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Microsoft Azure (Dynamic web app hosting on a tomcat server)</w:t></w:r></w:p><w:p><w:r><w:t>Tools for preparing presentation artefacts:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Adobe Photoshop</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Adobe Illustrator</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Adobe XD</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Visual Studio Code</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr></w:p>'

$ooxml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newContent + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$insertionRange.InsertXML($ooxml)
